$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct a student's email address (in place, same shared-string slot).
$ws.Range("B2").Value = "aswanibolisetti@gmail.com"

# Two more "registered students" rows get appended to the roster.
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "student@example.com"
$ws.Range("A8").Value = 7

# Make every data row's height explicit (matches a full Excel re-save).
$ws.Range("A1:B8").EntireRow.RowHeight = 15

# Leave the selection where the user's last edit was - cell B8.
$ws.Range("B8").Select()
